$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 656, shifting the existing rows 656-697 down
# to 657-698 (matches the diff: a new data point for 2026/01/19 16:00 was
# recorded, pushing everything from 2026/12/29 onward down by one row).
$ws.Rows.Item(656).Insert()

# Populate the newly inserted row 656 with its values.
# A leading apostrophe keeps this date-like string stored as text (like all
# the other date cells in column A) instead of Excel auto-converting it to a
# date serial number.
$ws.Range("A656").Value = "'2026/01/19"
$ws.Range("B656").Value = "月"
$ws.Range("C656").Value = 16
$ws.Range("D656").Value = 12
